# Updates the cryptocurrency price/volume table (cols D/E, rows 2-51) on
# Sheet1 to reflect refreshed market data from the scheduled GitHub Actions
# scraper run. Values are written as literal text (matching the workbook's
# existing inlineStr cell type) rather than numbers: most Price figures
# contain locale thousands-separator dots (e.g. "63.555.86") that must be
# preserved verbatim, and Volume(1h) values keep their padding spaces
# around the percentage. For Price cells whose new text otherwise looks
# like a plain number (e.g. "155.09", "37.00", "0.999") we briefly force a
# text NumberFormat so Excel stores them as text instead of auto-converting
# to a numeric value (which would also silently drop significant trailing
# zeros), then restore the cell's normal style so no formatting residue is
# left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.555.86"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "3.069.14"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").Value = "3.068.58"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000238"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "3.576.50"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "63.452.24"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "3.068.78"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "493.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.13%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("D36").Value = "0.0₃0826"
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "438.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  +2.63%  "
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "2.843.07"
$ws.Range("E46").Value = "  +1.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("E51").Value = "  -0.79%  "

